# Updating proposal list for the final subset and fixing hyperlinks.
#
# The old row 51 ("NOVEMBER Chair's text - sustainability criteria ...")
# is removed entirely; the row that used to be 52 ("Hide" helper row)
# shifts up and becomes the new row 51. Deleting the whole row (rather
# than just clearing it) is what causes Excel to drop the now-unused
# shared strings and re-pack the shared string table / re-number the
# sheet dimension, exactly as in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(51).Delete()

# Reflect the resulting scroll position/selection: after the edit the
# user ends up with C50 selected near the bottom of the (now 51-row)
# sheet.
$ws.Range("C50").Select()
